$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 46-48: coin rotates position (NEARProtocol -> WEMIXToken -> Maker -> NEARProtocol)
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.18"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.219.27"
$ws.Range("E47").Value = "  +4.99%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.57"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +5.85%  "

# Price (D) and Volume(1h) (E) updates for remaining rows
$ws.Range("D2").Value = "52.328.61"
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("D3").Value = "2.920.87"
$ws.Range("E3").Value = "  +4.36%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "354.14"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "112.99"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.71%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.562"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.49%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.632"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.42"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0867"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.43%  "
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.14"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.88"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.73%  "
$ws.Range("D15").Value = "3.385.52"
$ws.Range("E15").Value = "  +4.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.00"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +6.28%  "
$ws.Range("D17").Value = "2.915.92"
$ws.Range("E17").Value = "  +4.83%  "
$ws.Range("D18").Value = "52.390.54"
$ws.Range("E18").Value = "  +1.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.71"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.34"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +6.11%  "
$ws.Range("E21").Value = "  +7.25%  "
$ws.Range("D22").Value = "0.0₃0987"
$ws.Range("E22").Value = "  +1.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.22"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "271.89"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.57%  "
$ws.Range("E25").Value = "  +2.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.94"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.41%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("E28").Value = "  +1.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.66"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +3.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.15"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.90%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.56"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +6.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.25"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.18"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +7.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0951"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +10.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "53.12"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0457"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.35"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +7.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.98"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.09"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.84"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +15.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "23.97"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +9.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.117"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.65"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +8.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "121.47"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.262"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +24.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0340"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +14.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.972"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +8.20%  "
